$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$values = @{
    2  = "Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Majorelle Magdy, Dr. Servinaz Sayed Mohammad"
    3  = "Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy"
    4  = "Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy"
    5  = "Dr. Nourhan Mahmoud, Dr. Hanan Ragab, Dr. Eman Tantawi, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Nesma"
    6  = "Dr. Nourhan Mahmoud, Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Nahla Nagiub, Dr. Menna tuâ€™Allah Medhat, Dr. Asmaa Reda, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy"
    7  = "Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy"
    8  = "Dr. Shimaa Ahmad Mekki, Dr. Eman Tantawi, Dr. Manar Montaser, Administrator, Dr. Asmaa Reda, Dr. Majorelle Magdy"
    9  = "Dr. Manar Montaser, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Majorelle Magdy"
    10 = "Dr. Sara Wael, Dr. Shimaa Ahmad Mekki, Dr. Heba Mahmoud Ali, Dr. Alshimaa Atef, Dr. Gehan Adel, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad"
    11 = "Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Asmaa Reda"
    12 = "Dr. Salma El-Gendy, Administrator"
    13 = "Dr. Mariam Nour El-Din, D Wessam Atef, Dr. Shimaa Ashraf, Dr. Omnia Mohammad, Dr. Safa Hany"
    15 = "Dr. Amal Awwad, D Wessam Atef"
    17 = "Dr. Esraa Mostafa, Dr. Yasmeena Fattoh, Dr. Basma Hamed, Dr. Nourhan Osama, Dr. Madeha Saeed, Dr. Eman M. Abo-Sakaya, Dr. Marwa Mustafa, Dr. Dina Adel, Dr. Sarah Abdelmohsen, Dr. Arwa Al-Sayed"
    19 = "D Mariam E. Mohammad, Dr. Sarah Mahdy"
    22 = "Dr. Amr Saeed, Dr. Nancy Abd Al-Shafy"
    23 = "Dr. Nourham Mostafa, Dr. Hana Amr"
    24 = "Dr. Remon, Dr. Wafaa Ebida, Dr. Youstina Magdy, Dr. Salma Hassan, Dr. Aya Emad, Dr. Yasmin, Dr. Ola Abd Al-Fattah, Dr. Maryam Ashraf, Dr. Marina Atef, Dr. Neveen Nashaat, Dr. Monica"
    25 = "Dr. Remon, Dr. Youstina Magdy, Dr. Abdullah El-Agrody, Dr. Aya Emad, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Marina Atef"
    26 = "Dr. Youstina Magdy, Dr. Gehad Salah"
    27 = "Dr. Remon, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Eman Mohammad Al"
    28 = "Dr. Remon, Dr. Aya Hanafy, Dr. Wafaa Ebida, Dr. Abdullah El-Agrody, Dr. Salma Hassan, Dr. Eman Samir Gabry, Dr. Nardine, Dr. Neveen Nashaat"
    29 = "Dr. Remon, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Neveen Nashaat, Dr. Monica"
    30 = "Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy"
    31 = "Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy"
    32 = "Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy"
    33 = "Dr. Nourhan Mahmoud, Dr. Hanan Ragab, Dr. Eman Tantawi, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Nesma"
    34 = "Dr. Nourhan Mahmoud, Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Nahla Nagiub, Dr. Menna tuâ€™Allah Medhat, Dr. Asmaa Reda, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy"
    35 = "Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy"
    36 = "Dr. Shimaa Ahmad Mekki, Dr. Eman Tantawi, Dr. Manar Montaser, Administrator, Dr. Asmaa Reda, Dr. Majorelle Magdy"
    37 = "Dr. Manar Montaser, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Majorelle Magdy"
    38 = "Dr. Sara Wael, Dr. Shimaa Ahmad Mekki, Dr. Heba Mahmoud Ali, Dr. Alshimaa Atef, Dr. Gehan Adel, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad"
    39 = "Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Asmaa Reda"
    40 = "Dr. Salma El-Gendy, Administrator"
    41 = "Dr. Mariam Nour El-Din, D Wessam Atef, Dr. Shimaa Ashraf, Dr. Omnia Mohammad, Dr. Safa Hany"
    43 = "Dr. Amal Awwad, D Wessam Atef"
    45 = "Dr. Esraa Mostafa, Dr. Yasmeena Fattoh, Dr. Basma Hamed, Dr. Nourhan Osama, Dr. Madeha Saeed, Dr. Eman M. Abo-Sakaya, Dr. Marwa Mustafa, Dr. Dina Adel, Dr. Sarah Abdelmohsen, Dr. Arwa Al-Sayed"
    47 = "D Mariam E. Mohammad, Dr. Sarah Mahdy"
    50 = "Dr. Amr Saeed, Dr. Nancy Abd Al-Shafy"
    51 = "Dr. Nourham Mostafa, Dr. Hana Amr"
    52 = "Dr. Remon, Dr. Wafaa Ebida, Dr. Youstina Magdy, Dr. Salma Hassan, Dr. Aya Emad, Dr. Yasmin, Dr. Ola Abd Al-Fattah, Dr. Maryam Ashraf, Dr. Marina Atef, Dr. Neveen Nashaat, Dr. Monica"
    53 = "Dr. Remon, Dr. Youstina Magdy, Dr. Abdullah El-Agrody, Dr. Aya Emad, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Marina Atef"
    54 = "Dr. Youstina Magdy, Dr. Gehad Salah"
    55 = "Dr. Remon, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Eman Mohammad Al"
    56 = "Dr. Remon, Dr. Aya Hanafy, Dr. Wafaa Ebida, Dr. Abdullah El-Agrody, Dr. Salma Hassan, Dr. Eman Samir Gabry, Dr. Nardine, Dr. Neveen Nashaat"
    57 = "Dr. Remon, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Neveen Nashaat, Dr. Monica"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
